# Append the "Links:" section (github / demo hyperlinks) at the end of the
# document, after the existing "active should not allow null" paragraph.

$d = $word.ActiveDocument

$xml = @'
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">
        <w:body>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:sz w:val="18"/>
                <w:szCs w:val="18"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rtl w:val="0"/>
              </w:rPr>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:b w:val="1"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:u w:val="single"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rtl w:val="0"/>
              </w:rPr>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:b w:val="1"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:u w:val="single"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:b w:val="1"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:u w:val="single"/>
                <w:rtl w:val="0"/>
              </w:rPr>
              <w:t xml:space="preserve">Links:</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:sz w:val="18"/>
                <w:szCs w:val="18"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rtl w:val="0"/>
              </w:rPr>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:ind w:left="720" w:firstLine="0"/>
              <w:rPr>
                <w:color w:val="1155cc"/>
                <w:sz w:val="18"/>
                <w:szCs w:val="18"/>
                <w:highlight w:val="white"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:color w:val="24292f"/>
                <w:sz w:val="18"/>
                <w:szCs w:val="18"/>
                <w:highlight w:val="white"/>
                <w:rtl w:val="0"/>
              </w:rPr>
              <w:t xml:space="preserve">github: </w:t>
            </w:r>
            <w:hyperlink r:id="hlinkGithub">
              <w:r>
                <w:rPr>
                  <w:color w:val="1155cc"/>
                  <w:sz w:val="18"/>
                  <w:szCs w:val="18"/>
                  <w:highlight w:val="white"/>
                  <w:rtl w:val="0"/>
                </w:rPr>
                <w:t xml:space="preserve">https://github.com/gh11345/dlc</w:t>
              </w:r>
            </w:hyperlink>
            <w:r>
              <w:rPr>
                <w:rtl w:val="0"/>
              </w:rPr>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:ind w:left="720" w:firstLine="0"/>
              <w:rPr>
                <w:sz w:val="12"/>
                <w:szCs w:val="12"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:color w:val="24292f"/>
                <w:sz w:val="18"/>
                <w:szCs w:val="18"/>
                <w:highlight w:val="white"/>
                <w:rtl w:val="0"/>
              </w:rPr>
              <w:t xml:space="preserve">demo: </w:t>
            </w:r>
            <w:hyperlink r:id="hlinkDemo">
              <w:r>
                <w:rPr>
                  <w:color w:val="1155cc"/>
                  <w:sz w:val="18"/>
                  <w:szCs w:val="18"/>
                  <w:highlight w:val="white"/>
                  <w:rtl w:val="0"/>
                </w:rPr>
                <w:t xml:space="preserve">http://34.212.176.88/</w:t>
              </w:r>
            </w:hyperlink>
            <w:r>
              <w:rPr>
                <w:rtl w:val="0"/>
              </w:rPr>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
  <pkg:part pkg:name="/word/_rels/document.xml.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml">
    <pkg:xmlData>
      <Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships">
        <Relationship Id="hlinkGithub" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="https://github.com/gh11345/dlc" TargetMode="External"/>
        <Relationship Id="hlinkDemo" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="http://34.212.176.88/" TargetMode="External"/>
      </Relationships>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$r = $d.Content
$r.Collapse(0)
$r.InsertXML($xml)

# The XML round-trip that InsertXML performs treats an explicit
# w:firstLine="0" as the (omitted) default, so it gets dropped from the
# two "github:"/"demo:" link paragraphs. Restore it through the
# ParagraphFormat OM so the indentation markup matches the source exactly.
$total = $d.Paragraphs.Count
$githubPara = $d.Paragraphs.Item($total - 1)
$demoPara = $d.Paragraphs.Item($total)
$githubPara.Range.ParagraphFormat.FirstLineIndent = 0
$demoPara.Range.ParagraphFormat.FirstLineIndent = 0
